{"js": "// This document is a daily \"addition & subtraction within 100\" drill sheet:\n// one heading paragraph with the date, followed by a single 20x5 table of\n// arithmetic problems. The edit (1) bumps the heading to the next day and\n// (2) replaces every problem in the table with a new one, cell-for-cell,\n// while leaving all existing fonts/sizes/alignment untouched.\n\nconst body = context.document.body;\n\n// --- 1. Update the date heading (first paragraph of the body) -------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2024-12-09 Monday\", Word.InsertLocation.replace);\n\n// --- 2. Update every cell of the practice-problem table --------------------\n// Setting Table.values rewrites each cell's text in place (row-major,\n// left-to-right, top-to-bottom) without touching the cell/run formatting\n// already present in the document.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"2+40=\", \"91-11=\", \"19+54=\", \"7-0=\", \"69-57=\"],\n  [\"20+3=\", \"12-4=\", \"40-28=\", \"26+1=\", \"55-32=\"],\n  [\"54-1=\", \"30+44=\", \"70+16=\", \"17+57=\", \"62+28=\"],\n  [\"94-15=\", \"72-67=\", \"28+44=\", \"59-44=\", \"86-61=\"],\n  [\"52+22=\", \"51-38=\", \"86-12=\", \"7+3=\", \"42+38=\"],\n  [\"92-72=\", \"42-1=\", \"23+35=\", \"27+42=\", \"14+80=\"],\n  [\"78-31=\", \"93-72=\", \"29+14=\", \"36-18=\", \"53-12=\"],\n  [\"2+23=\", \"61-25=\", \"41+0=\", \"54-33=\", \"52-8=\"],\n  [\"85-54=\", \"39+3=\", \"13+23=\", \"22+0=\", \"64-0=\"],\n  [\"8+11=\", \"36+50=\", \"40+49=\", \"58+9=\", \"60+18=\"],\n  [\"16-11=\", \"56-37=\", \"36-1=\", \"56-33=\", \"34-7=\"],\n  [\"10+13=\", \"74-14=\", \"70-29=\", \"36+14=\", \"10+56=\"],\n  [\"78-59=\", \"56-11=\", \"17-16=\", \"80-64=\", \"73-6=\"],\n  [\"94-52=\", \"10+58=\", \"29+60=\", \"43+54=\", \"2+40=\"],\n  [\"43-1=\", \"39-18=\", \"63-36=\", \"3-0=\", \"85-65=\"],\n  [\"46+10=\", \"30+57=\", \"80-6=\", \"28+57=\", \"1+35=\"],\n  [\"47-5=\", \"23+57=\", \"84+11=\", \"27+48=\", \"8+0=\"],\n  [\"64-9=\", \"10+54=\", \"28-12=\", \"38+2=\", \"44+34=\"],\n  [\"97-3=\", \"32+46=\", \"33+27=\", \"84-39=\", \"39+31=\"],\n  [\"13+27=\", \"9+80=\", \"30+42=\", \"33-21=\", \"3+37=\"]\n];\n\nawait context.sync();\n", "ps1": "# This document is a daily \"addition & subtraction within 100\" drill sheet:\n# one heading paragraph with the date, followed by a single 20x5 table of\n# arithmetic problems. The edit (1) bumps the heading to the next day and\n# (2) replaces every problem in the table with a new one, cell-for-cell,\n# while leaving all existing fonts/sizes/alignment untouched.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date heading (first paragraph of the body) --------------\n# Assigning Range.Text replaces only the paragraph's text, leaving the\n# paragraph mark (and its run formatting) in place.\n$d.Paragraphs.Item(1).Range.Text = \"2024-12-09 Monday\"\n\n# --- 2. Update every cell of the practice-problem table --------------------\n# Assigning each Cell's Range.Text rewrites only the cell's text, leaving\n# the cell/run formatting already present in the document untouched.\n$t = $d.Tables.Item(1)\n$values = @(\n    @(\"2+40=\", \"91-11=\", \"19+54=\", \"7-0=\", \"69-57=\"),\n    @(\"20+3=\", \"12-4=\", \"40-28=\", \"26+1=\", \"55-32=\"),\n    @(\"54-1=\", \"30+44=\", \"70+16=\", \"17+57=\", \"62+28=\"),\n    @(\"94-15=\", \"72-67=\", \"28+44=\", \"59-44=\", \"86-61=\"),\n    @(\"52+22=\", \"51-38=\", \"86-12=\", \"7+3=\", \"42+38=\"),\n    @(\"92-72=\", \"42-1=\", \"23+35=\", \"27+42=\", \"14+80=\"),\n    @(\"78-31=\", \"93-72=\", \"29+14=\", \"36-18=\", \"53-12=\"),\n    @(\"2+23=\", \"61-25=\", \"41+0=\", \"54-33=\", \"52-8=\"),\n    @(\"85-54=\", \"39+3=\", \"13+23=\", \"22+0=\", \"64-0=\"),\n    @(\"8+11=\", \"36+50=\", \"40+49=\", \"58+9=\", \"60+18=\"),\n    @(\"16-11=\", \"56-37=\", \"36-1=\", \"56-33=\", \"34-7=\"),\n    @(\"10+13=\", \"74-14=\", \"70-29=\", \"36+14=\", \"10+56=\"),\n    @(\"78-59=\", \"56-11=\", \"17-16=\", \"80-64=\", \"73-6=\"),\n    @(\"94-52=\", \"10+58=\", \"29+60=\", \"43+54=\", \"2+40=\"),\n    @(\"43-1=\", \"39-18=\", \"63-36=\", \"3-0=\", \"85-65=\"),\n    @(\"46+10=\", \"30+57=\", \"80-6=\", \"28+57=\", \"1+35=\"),\n    @(\"47-5=\", \"23+57=\", \"84+11=\", \"27+48=\", \"8+0=\"),\n    @(\"64-9=\", \"10+54=\", \"28-12=\", \"38+2=\", \"44+34=\"),\n    @(\"97-3=\", \"32+46=\", \"33+27=\", \"84-39=\", \"39+31=\"),\n    @(\"13+27=\", \"9+80=\", \"30+42=\", \"33-21=\", \"3+37=\")\n)\n\nfor ($r = 0; $r -lt $values.Count; $r++) {\n    $row = $values[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
